$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 16.905
$ws.Range("C3").Value = -13.666
$ws.Range("E3").Value = 16.147
$ws.Range("C4").Value = -12.168
$ws.Range("E9").Value = 17.306
$ws.Range("B11").Value = 6.825
$ws.Range("B12").Value = 4.872
$ws.Range("C14").Value = -11.939
$ws.Range("B15").Value = 5.340000000000001
$ws.Range("E15").Value = 16.13
$ws.Range("E19").Value = 16.35
$ws.Range("E20").Value = 16.417
$ws.Range("E25").Value = 17.134
$ws.Range("C26").Value = -11.922
$ws.Range("B27").Value = 5.386
$ws.Range("E27").Value = 16.056
$ws.Range("B28").Value = 5.537
$ws.Range("E28").Value = 16.957
$ws.Range("E30").Value = 16.033
$ws.Range("B31").Value = 5.099000000000001
$ws.Range("C31").Value = -12.819
$ws.Range("B32").Value = 7.645999999999999
$ws.Range("E32").Value = 16.595
$ws.Range("C35").Value = -12.654
$ws.Range("B36").Value = 8.888
$ws.Range("C37").Value = -13.465
$ws.Range("B38").Value = 5.514
$ws.Range("C39").Value = -12.746
$ws.Range("C40").Value = -12.51
$ws.Range("E44").Value = 16.749
$ws.Range("C45").Value = -12.651
$ws.Range("B46").Value = 6.540999999999999
$ws.Range("E47").Value = 16.224
$ws.Range("C52").Value = -11.066
$ws.Range("B54").Value = 5.118
$ws.Range("B55").Value = 4.74
$ws.Range("B56").Value = 4.249000000000001
$ws.Range("C57").Value = -13.697
$ws.Range("E58").Value = 16.607
$ws.Range("E62").Value = 16.4
$ws.Range("B67").Value = 5.331999999999999
$ws.Range("B69").Value = 5.332
$ws.Range("B72").Value = 5.616
$ws.Range("B73").Value = 7.582000000000001
$ws.Range("E77").Value = 17.045
$ws.Range("E78").Value = 16.394
$ws.Range("C81").Value = -13.358
$ws.Range("B83").Value = 5.290999999999999
$ws.Range("C83").Value = -14.083
$ws.Range("E84").Value = 16.384
$ws.Range("B86").Value = 4.94
$ws.Range("E89").Value = 17.158
$ws.Range("B91").Value = 5.581
$ws.Range("E91").Value = 17.429
$ws.Range("E92").Value = 17.373
$ws.Range("B93").Value = 5.621
$ws.Range("E96").Value = 16.512
$ws.Range("B99").Value = 5.220000000000001
$ws.Range("C100").Value = -12.691
$ws.Range("C102").Value = -12.637
$ws.Range("E102").Value = 16.523
